$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "TreatmentTab" row (row 5) query in column B wrapped REPLACE() in a
# redundant CONCAT() call. Remove the unnecessary CONCAT() wrapper.
$oldQuery = $ws.Range("B5").Value2
$newQuery = $oldQuery.Replace("CONCAT(REPLACE(trt.treatment_agent, ';', ', '))", "REPLACE(trt.treatment_agent, ';', ', ')")
$ws.Range("B5").Value2 = $newQuery

# Re-apply the (visually unchanged) font/wrap formatting so the cell picks
# up a freshly written style record, matching how Excel re-stamps a cell's
# format after its content is edited.
$ws.Range("B5").Font.ThemeColor = 1
$ws.Range("B5").Font.Size = 12
$ws.Range("B5").WrapText = $true

# Move the active selection to B2 (the sheet's top-left cell is no longer
# scrolled/frozen away from A1).
$null = $ws.Range("B2").Select()
